$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.524.64'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.472.32'
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.78'
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.91'
$ws.Range("E6").Value = '  -3.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").Value = '  +2.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.47'
$ws.Range("E10").Value = '  -3.51%  '

$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("E12").Value = '  +1.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.854.66'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.86'
$ws.Range("E14").Value = '  -2.14%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.18'
$ws.Range("E15").Value = '  +6.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.482.04'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.767'
$ws.Range("E17").Value = '  -3.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.511.56'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.50'
$ws.Range("E19").Value = '  +2.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("E20").Value = '  +2.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.03'
$ws.Range("E21").Value = '  +4.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.08'
$ws.Range("E22").Value = '  -1.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.80'
$ws.Range("E23").Value = '  -0.86%  '

$ws.Range("E24").Value = '  -1.26%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("E26").Value = '  -1.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.84'
$ws.Range("E27").Value = '  +2.79%  '

$ws.Range("E28").Value = '  -0.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.69'
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.46'
$ws.Range("E30").Value = '  -3.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.32'
$ws.Range("E31").Value = '  +2.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.45'
$ws.Range("E32").Value = '  -0.87%  '

$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0758'
$ws.Range("E34").Value = '  +1.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.18'
$ws.Range("E35").Value = '  -3.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.35'
$ws.Range("E36").Value = '  -8.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.88'
$ws.Range("E37").Value = '  -6.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("E38").Value = '  +1.17%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.114'
$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.80'
$ws.Range("E40").Value = '  -4.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  -6.25%  '

$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.958.83'
$ws.Range("E43").Value = '  -1.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0284'
$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.61'
$ws.Range("E45").Value = '  -5.50%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.93'
$ws.Range("E46").Value = '  -2.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.07'
$ws.Range("E47").Value = '  +3.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.713.69'
$ws.Range("E48").Value = '  -0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.33'
$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.31'
$ws.Range("E50").Value = '  -3.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.02'
$ws.Range("E51").Value = '  -3.79%  '
